# "Generate Report for Handoff"
#
# The row for file dbe0b373-f1cd-4dad-a6c4-166ceca49656 (row 7 on every
# sheet) gets fresh handoff timestamps stamped onto it:
#   - Overview!G7            "Latest HO Xliff Generate Date"
#   - zh-cn!H7                "Latest Handoff Datetime"
#   - de-de!H7                "Latest Handoff Datetime"
#
# These columns are formatted as plain text (not real date serials), so the
# values below must land as text, matching the existing "yyyy-mm-dd HH:mm:ss"-
# look-alike strings already in the sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-09-06 10:55:55"

$zhCn = $wb.Worksheets.Item("zh-cn")
$zhCn.Range("H7").Value = "2016-09-06 10:55:50"

$deDe = $wb.Worksheets.Item("de-de")
$deDe.Range("H7").Value = "2016-09-06 10:55:55"
